$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (2-7) with new Documento/Grupo values ---
$ws.Range("A2").Value = 1023866054
$ws.Range("B2").Value = 196

$ws.Range("A3").Value = 1052412499
$ws.Range("B3").Value = 196

$ws.Range("A4").Value = 1098658561
$ws.Range("B4").Value = 195

$ws.Range("A5").Value = 63528540
$ws.Range("B5").Value = 195

$ws.Range("A6").Value = 39427884
$ws.Range("B6").Value = 194

$ws.Range("A7").Value = 1040370636
$ws.Range("B7").Value = 194

# --- Row 8: new values, Perfil bumps from 3 to 4 ---
$ws.Range("A8").Value = 80125674
$ws.Range("B8").Value = 196
$ws.Range("C8").Value = 4

# --- New rows 9-11 ---
$ws.Range("A9").Value = 1032445017
$ws.Range("B9").Value = 196
$ws.Range("C9").Value = 4

$ws.Range("A10").Value = 1014186124
$ws.Range("B10").Value = 195
$ws.Range("C10").Value = 4

$ws.Range("A11").Value = 39424174
$ws.Range("B11").Value = 194
$ws.Range("C11").Value = 4

# --- Re-apply the "Normal" style on rows 4-8 (Documento/Perfil columns, plus
#     Grupo on row 8) so they pick up a distinct cell format record, matching
#     the source workbook's re-saved style table. ---
$ws.Range("A4:A8").Style = "Normal"
$ws.Range("C4:C8").Style = "Normal"
$ws.Range("B8").Style = "Normal"

# --- Selection ends on B11, matching the last-edited cell ---
$ws.Range("B11").Select() | Out-Null
